$wb = $excel.ActiveWorkbook

$oldGuid = "432182e3-a89b-4f3d-8ec0-3c915ef01741"
$newGuid = "390b7c5e-abd0-4eb9-a7de-e37cf731c959"

$oldHash = "2acf353bf28ad25d265d96beba17465d5ea6c129"
$newHash = "2182be7151eae1ad44bf34891b67bf8bdcee4459"

$newHandoffDate   = "2016-08-31 12:41:40"
$newZhHandoffDate = "2016-08-31 12:41:28"

$hyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/13f2a2d7cb94eef4bad58039efdb94e71d95d7ef/e2e/$oldGuid.md"

# --- Sheet "Overview" ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = $newHandoffDate

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hyperlinkAddress, "", "", "e2e\$newGuid.md")

# --- Sheet "zh-cn" ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = "$newGuid.md"
$wsZh.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZh.Range("H2").Value = $newZhHandoffDate

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $hyperlinkAddress, "", "", "$newGuid.md")

# --- Sheet "de-de" ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = "$newGuid.md"
$wsDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDe.Range("H2").Value = $newHandoffDate

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $hyperlinkAddress, "", "", "$newGuid.md")
